$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 316.33334
$ws.Range("I12").Value = 99
$ws.Range("J12").Value = 425
$ws.Range("K12").Value = 99
$ws.Range("L12").Value = 425
$ws.Range("M12").Value = 71
$ws.Range("N12").Value = -765
$ws.Range("H46").Value = 24475.117
$ws.Range("I46").Value = 817
$ws.Range("J46").Value = 25953.75
$ws.Range("K46").Value = 2451
$ws.Range("L46").Value = 77861.25
$ws.Range("M46").Value = -2332
$ws.Range("N46").Value = -78099.25
$ws.Range("H60").Value = 24475.117
$ws.Range("I60").Value = 817
$ws.Range("J60").Value = 25953.75
$ws.Range("K60").Value = 2451
$ws.Range("L60").Value = 77861.25
$ws.Range("M60").Value = -1967
$ws.Range("N60").Value = -78829.25
$ws.Range("H87").Value = 15513.274
$ws.Range("J87").Value = 15513.274
$ws.Range("L87").Value = 15513.274
$ws.Range("N87").Value = -18009.274
$ws.Range("H90").Value = 15513.274
$ws.Range("J90").Value = 15513.274
$ws.Range("L90").Value = 46539.822
$ws.Range("N90").Value = -59019.822
$ws.Range("H112").Value = 1309.4
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 1273.0952
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 3819.2856
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -6035.2856
$ws.Range("H113").Value = 7196.6562
$ws.Range("I113").Value = 2752.9333
$ws.Range("J113").Value = 11117.588
$ws.Range("K113").Value = 2752.9333
$ws.Range("L113").Value = 11117.588
$ws.Range("M113").Value = 501.0666999999999
$ws.Range("N113").Value = -17625.588
$ws.Range("H129").Value = 848.0270400000001
$ws.Range("J129").Value = 888.74243
$ws.Range("L129").Value = 2666.22729
$ws.Range("N129").Value = -12666.22729
$ws.Range("H132").Value = 37240.37
$ws.Range("I132").Value = 58169.047
$ws.Range("J132").Value = 1822.6154
$ws.Range("K132").Value = 174507.141
$ws.Range("L132").Value = 5467.8462
$ws.Range("M132").Value = -171977.141
$ws.Range("N132").Value = -10527.8462
$ws.Range("H138").Value = 3220.87
$ws.Range("J138").Value = 4781.469
$ws.Range("L138").Value = 14344.407
$ws.Range("N138").Value = -24624.407

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 47501.5
$ws.Range("H34").Value = 25000
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H61").Value = 3892.926
$ws.Range("I61").Value = 4372.722
$ws.Range("J61").Value = 2933.3333
$ws.Range("K61").Value = 4372.722
$ws.Range("L61").Value = 2933.3333
$ws.Range("M61").Value = -4160.722
$ws.Range("N61").Value = -3357.3333
$ws.Range("H132").Value = 1535.5781
$ws.Range("I132").Value = 1203.3778
$ws.Range("J132").Value = 2322.3684
$ws.Range("K132").Value = 3610.1334
$ws.Range("L132").Value = 6967.1052
$ws.Range("M132").Value = -1080.1334
$ws.Range("N132").Value = -12027.1052
$ws.Range("H136").Value = 3892.926
$ws.Range("I136").Value = 4372.722
$ws.Range("J136").Value = 2933.3333
$ws.Range("K136").Value = 13118.166
$ws.Range("L136").Value = 8799.999899999999
$ws.Range("M136").Value = -10568.166
$ws.Range("N136").Value = -13899.9999
$ws.Range("H139").Value = 49283.332
$ws.Range("J139").Value = 49283.332
$ws.Range("L139").Value = 49283.332
$ws.Range("N139").Value = -59563.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 2666.6667
$ws.Range("I14").Value = 2000
$ws.Range("J14").Value = 4000
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = -1828
$ws.Range("N14").Value = -4344
$ws.Range("H80").Value = 2929.037
$ws.Range("I80").Value = 569.6
$ws.Range("J80").Value = 4316.9414
$ws.Range("K80").Value = 569.6
$ws.Range("L80").Value = 4316.9414
$ws.Range("M80").Value = 428.4
$ws.Range("N80").Value = -6312.9414
$ws.Range("H83").Value = 2929.037
$ws.Range("I83").Value = 569.6
$ws.Range("J83").Value = 4316.9414
$ws.Range("K83").Value = 2848
$ws.Range("L83").Value = 21584.707
$ws.Range("M83").Value = 2144
$ws.Range("N83").Value = -31568.707
$ws.Range("H107").Value = 880.2857
$ws.Range("I107").Value = 860.3333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 860.3333
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1059.6667
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 1165.6904
$ws.Range("I134").Value = 904.53125
$ws.Range("J134").Value = 2001.4
$ws.Range("K134").Value = 2713.59375
$ws.Range("L134").Value = 6004.200000000001
$ws.Range("M134").Value = -178.59375
$ws.Range("N134").Value = -11074.2
$ws.Range("H140").Value = 54830
$ws.Range("J140").Value = 54830
$ws.Range("L140").Value = 54830
$ws.Range("N140").Value = -65190

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1658.4286
$ws.Range("I22").Value = 2221.8
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 2221.8
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -1871.8
$ws.Range("N22").Value = -950
$ws.Range("H35").Value = 1150
$ws.Range("I35").Value = 1150
$ws.Range("K35").Value = 1150
$ws.Range("M35").Value = -856
$ws.Range("H107").Value = 384.2857
$ws.Range("I107").Value = 295.375
$ws.Range("J107").Value = 502.83334
$ws.Range("K107").Value = 295.375
$ws.Range("L107").Value = 502.83334
$ws.Range("M107").Value = 1624.625
$ws.Range("N107").Value = -4342.83334
$ws.Range("H134").Value = 1750.421
$ws.Range("I134").Value = 1729.0286
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 5187.085800000001
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -2652.085800000001
$ws.Range("N134").Value = -11070
$ws.Range("H140").Value = 67023.37
$ws.Range("J140").Value = 67023.37
$ws.Range("L140").Value = 67023.37
$ws.Range("N140").Value = -77383.37

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 32.5
$ws.Range("K7").Value = 97.5
$ws.Range("M7").Value = 14.5
$ws.Range("H23").Value = 89.583336
$ws.Range("J23").Value = 97
$ws.Range("L23").Value = 291
$ws.Range("N23").Value = -761
$ws.Range("H34").Value = 670.8889
$ws.Range("J34").Value = 759.7273
$ws.Range("L34").Value = 2279.1819
$ws.Range("N34").Value = -2447.1819
$ws.Range("H113").Value = 590.2
$ws.Range("I113").Value = 494.25
$ws.Range("J113").Value = 699.8570999999999
$ws.Range("K113").Value = 1482.75
$ws.Range("L113").Value = 2099.5713
$ws.Range("M113").Value = 687.25
$ws.Range("N113").Value = -6439.5713
$ws.Range("H131").Value = 877.55
$ws.Range("J131").Value = 883.3838500000001
$ws.Range("L131").Value = 2650.15155
$ws.Range("N131").Value = -12730.15155
$ws.Range("H132").Value = 2030.2727
$ws.Range("I132").Value = 1862.4
$ws.Range("J132").Value = 2079.647
$ws.Range("K132").Value = 16761.6
$ws.Range("L132").Value = 18716.823
$ws.Range("M132").Value = -14231.6
$ws.Range("N132").Value = -23776.823
$ws.Range("H134").Value = 1584.48
$ws.Range("I134").Value = 1256.2222
$ws.Range("J134").Value = 2428.5715
$ws.Range("K134").Value = 3768.6666
$ws.Range("L134").Value = 7285.7145
$ws.Range("M134").Value = 1301.3334
$ws.Range("N134").Value = -17425.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5321.1567
$ws.Range("I70").Value = 4754.048
$ws.Range("J70").Value = 5718.1333
$ws.Range("K70").Value = 4754.048
$ws.Range("L70").Value = 5718.1333
$ws.Range("M70").Value = -4484.048
$ws.Range("N70").Value = -6258.1333
$ws.Range("H73").Value = 5321.1567
$ws.Range("I73").Value = 4754.048
$ws.Range("J73").Value = 5718.1333
$ws.Range("K73").Value = 4754.048
$ws.Range("L73").Value = 5718.1333
$ws.Range("M73").Value = -3818.048
$ws.Range("N73").Value = -7590.1333
$ws.Range("H132").Value = 2094.054
$ws.Range("I132").Value = 1732.8462
$ws.Range("J132").Value = 2947.818
$ws.Range("K132").Value = 5198.5386
$ws.Range("L132").Value = 8843.454000000002
$ws.Range("M132").Value = -2668.5386
$ws.Range("N132").Value = -13903.454
$ws.Range("H135").Value = 41176
$ws.Range("J135").Value = 41176
$ws.Range("L135").Value = 41176
$ws.Range("N135").Value = -51316
$ws.Range("H138").Value = 45994
$ws.Range("J138").Value = 45994
$ws.Range("L138").Value = 45994
$ws.Range("N138").Value = -56274

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 40041
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H61").Value = 253094.5
$ws.Range("I61").Value = 253094.5
$ws.Range("K61").Value = 253094.5
$ws.Range("M61").Value = -252892.5
$ws.Range("H113").Value = 253094.5
$ws.Range("I113").Value = 253094.5
$ws.Range("K113").Value = 253094.5
$ws.Range("M113").Value = -250924.5
$ws.Range("H127").Value = 33992.5
$ws.Range("J127").Value = 33992.5
$ws.Range("L127").Value = 33992.5
$ws.Range("N127").Value = -43912.5
$ws.Range("H132").Value = 6533.5674
$ws.Range("I132").Value = 6853.977
$ws.Range("J132").Value = 5920.609
$ws.Range("K132").Value = 20561.931
$ws.Range("L132").Value = 17761.827
$ws.Range("M132").Value = -18031.931
$ws.Range("N132").Value = -22821.827
$ws.Range("H133").Value = 43409.25
$ws.Range("J133").Value = 43409.25
$ws.Range("L133").Value = 43409.25
$ws.Range("N133").Value = -48469.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 16000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 16000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 16000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -16826
$ws.Range("H43").Value = 17333.334
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 17333.334
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 17333.334
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -17631.334
$ws.Range("H114").Value = 33718.4
$ws.Range("J114").Value = 33718.4
$ws.Range("L114").Value = 33718.4
$ws.Range("N114").Value = -42396.4
$ws.Range("H137").Value = 55148.832
$ws.Range("J137").Value = 55148.832
$ws.Range("L137").Value = 55148.832
$ws.Range("N137").Value = -65348.832
